$wb = $excel.ActiveWorkbook

# Sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H82").Value = 2439.1667
$ws.Range("I82").Value = 633.75
$ws.Range("J82").Value = 6050
$ws.Range("K82").Value = 1901.25
$ws.Range("L82").Value = 18150
$ws.Range("M82").Value = -1495.25
$ws.Range("N82").Value = -18962
$ws.Range("H85").Value = 2439.1667
$ws.Range("I85").Value = 633.75
$ws.Range("J85").Value = 6050
$ws.Range("K85").Value = 1901.25
$ws.Range("L85").Value = 18150
$ws.Range("M85").Value = -497.25
$ws.Range("N85").Value = -20958

# Sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5579.041
$ws.Range("I32").Value = 4048.054
$ws.Range("J32").Value = 10299.583
$ws.Range("K32").Value = 4048.054
$ws.Range("L32").Value = 10299.583
$ws.Range("M32").Value = -3761.054
$ws.Range("N32").Value = -10873.583
$ws.Range("H45").Value = 1736.2941
$ws.Range("I45").Value = 1081.2307
$ws.Range("J45").Value = 3865.25
$ws.Range("K45").Value = 1081.2307
$ws.Range("L45").Value = 3865.25
$ws.Range("M45").Value = -704.2307000000001
$ws.Range("N45").Value = -4619.25
$ws.Range("H61").Value = 4436.9375
$ws.Range("I61").Value = 1215.1666
$ws.Range("J61").Value = 6370
$ws.Range("K61").Value = 1215.1666
$ws.Range("L61").Value = 6370
$ws.Range("M61").Value = -1003.1666
$ws.Range("N61").Value = -6794
$ws.Range("H122").Value = 4407.857
$ws.Range("I122").Value = 2816
$ws.Range("J122").Value = 5999.7144
$ws.Range("K122").Value = 8448
$ws.Range("L122").Value = 17999.1432
$ws.Range("M122").Value = -5998
$ws.Range("N122").Value = -22899.1432
$ws.Range("H133").Value = 27980.555
$ws.Range("J133").Value = 27980.555
$ws.Range("L133").Value = 27980.555
$ws.Range("N133").Value = -33040.555
$ws.Range("H136").Value = 4436.9375
$ws.Range("I136").Value = 1215.1666
$ws.Range("J136").Value = 6370
$ws.Range("K136").Value = 3645.4998
$ws.Range("L136").Value = 19110
$ws.Range("M136").Value = -1095.4998
$ws.Range("N136").Value = -24210

# Sheet BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 823560.6
$ws.Range("I86").Value = 1151591.1
$ws.Range("J86").Value = 3484.5
$ws.Range("K86").Value = 1151591.1
$ws.Range("L86").Value = 3484.5
$ws.Range("M86").Value = -1150468.1
$ws.Range("N86").Value = -5730.5
$ws.Range("H89").Value = 823560.6
$ws.Range("I89").Value = 1151591.1
$ws.Range("J89").Value = 3484.5
$ws.Range("K89").Value = 5757955.5
$ws.Range("L89").Value = 17422.5
$ws.Range("M89").Value = -5752339.5
$ws.Range("N89").Value = -28654.5
$ws.Range("H105").Value = 1736.2
$ws.Range("I105").Value = 1544.375
$ws.Range("J105").Value = 2503.5
$ws.Range("K105").Value = 1544.375
$ws.Range("L105").Value = 2503.5
$ws.Range("M105").Value = 202.625
$ws.Range("N105").Value = -5997.5
$ws.Range("H134").Value = 2153.074
$ws.Range("I134").Value = 1060.591
$ws.Range("J134").Value = 6960
$ws.Range("K134").Value = 3181.773
$ws.Range("L134").Value = 20880
$ws.Range("M134").Value = -646.7729999999997
$ws.Range("N134").Value = -25950

# Sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2139.3408
$ws.Range("I31").Value = 1370.1111
$ws.Range("J31").Value = 3361.0588
$ws.Range("K31").Value = 1370.1111
$ws.Range("L31").Value = 3361.0588
$ws.Range("M31").Value = -1075.1111
$ws.Range("N31").Value = -3951.0588
$ws.Range("H34").Value = 2139.3408
$ws.Range("I34").Value = 1370.1111
$ws.Range("J34").Value = 3361.0588
$ws.Range("K34").Value = 1370.1111
$ws.Range("L34").Value = 3361.0588
$ws.Range("M34").Value = -1168.1111
$ws.Range("N34").Value = -3765.0588

# Sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 1965.1818
$ws.Range("I97").Value = 746.5
$ws.Range("J97").Value = 2236
$ws.Range("K97").Value = 2239.5
$ws.Range("L97").Value = 6708
$ws.Range("M97").Value = -1743.5
$ws.Range("N97").Value = -7700
$ws.Range("H107").Value = 1020.7273
$ws.Range("I107").Value = 740.55554
$ws.Range("J107").Value = 1092.7715
$ws.Range("K107").Value = 2221.66662
$ws.Range("L107").Value = 3278.3145
$ws.Range("M107").Value = -301.66662
$ws.Range("N107").Value = -7118.3145
$ws.Range("H113").Value = 981
$ws.Range("I113").Value = 549.25
$ws.Range("J113").Value = 1138
$ws.Range("K113").Value = 1647.75
$ws.Range("L113").Value = 3414
$ws.Range("M113").Value = 522.25
$ws.Range("N113").Value = -7754

# Sheet GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H95").Value = 16500
$ws.Range("J95").Value = 16500
$ws.Range("L95").Value = 16500
$ws.Range("N95").Value = -21992
$ws.Range("H102").Value = 2830.8235
$ws.Range("I102").Value = 2032.4
$ws.Range("J102").Value = 3971.4285
$ws.Range("K102").Value = 2032.4
$ws.Range("L102").Value = 3971.4285
$ws.Range("M102").Value = -410.4000000000001
$ws.Range("N102").Value = -7215.4285
$ws.Range("H132").Value = 23812738
$ws.Range("I132").Value = 41668336
$ws.Range("J132").Value = 5272.778
$ws.Range("K132").Value = 125005008
$ws.Range("L132").Value = 15818.334
$ws.Range("M132").Value = -125002478
$ws.Range("N132").Value = -20878.334

# Sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2832.5652
$ws.Range("I16").Value = 1649.9333
$ws.Range("J16").Value = 5050
$ws.Range("K16").Value = 1649.9333
$ws.Range("L16").Value = 5050
$ws.Range("M16").Value = -1479.9333
$ws.Range("N16").Value = -5390
$ws.Range("H22").Value = 1673.75
$ws.Range("I22").Value = 397.5
$ws.Range("J22").Value = 2950
$ws.Range("K22").Value = 397.5
$ws.Range("L22").Value = 2950
$ws.Range("M22").Value = -102.5
$ws.Range("N22").Value = -3540
$ws.Range("H27").Value = 1673.75
$ws.Range("I27").Value = 397.5
$ws.Range("J27").Value = 2950
$ws.Range("K27").Value = 397.5
$ws.Range("L27").Value = 2950
$ws.Range("M27").Value = -290.5
$ws.Range("N27").Value = -3164
$ws.Range("H46").Value = 1484.7142
$ws.Range("I46").Value = 820.7
$ws.Range("J46").Value = 3144.75
$ws.Range("K46").Value = 820.7
$ws.Range("L46").Value = 3144.75
$ws.Range("M46").Value = -632.7
$ws.Range("N46").Value = -3520.75
$ws.Range("H82").Value = 2663.611
$ws.Range("I82").Value = 2214.2
$ws.Range("J82").Value = 3225.375
$ws.Range("K82").Value = 2214.2
$ws.Range("L82").Value = 3225.375
$ws.Range("M82").Value = -1853.2
$ws.Range("N82").Value = -3947.375
$ws.Range("H85").Value = 2663.611
$ws.Range("I85").Value = 2214.2
$ws.Range("J85").Value = 3225.375
$ws.Range("K85").Value = 2214.2
$ws.Range("L85").Value = 3225.375
$ws.Range("M85").Value = -966.1999999999998
$ws.Range("N85").Value = -5721.375
$ws.Range("H136").Value = 2074.8333
$ws.Range("I136").Value = 1603.1875
$ws.Range("J136").Value = 3018.125
$ws.Range("K136").Value = 4809.5625
$ws.Range("L136").Value = 9054.375
$ws.Range("M136").Value = -2259.5625
$ws.Range("N136").Value = -14154.375

# Sheet WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").Value = -40040
